$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.600.41"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = "'1.854.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = "'265.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.02%  '

$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("E7").Value = '  -0.51%  '

$ws.Range("D8").Value = "'0.3290"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").Value = "'0.06815"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").Value = "'18.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.08%  '

$ws.Range("D11").Value = "'0.7777"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").Value = "'0.07769"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("D13").Value = "'1.850.67"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.60%  '

$ws.Range("D14").Value = "'88.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").Value = "'5.031"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.73%  '

$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = "'14.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("D18").Value = "'0.000007978"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").Value = "'0.9990"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").Value = "'26.614.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("D21").Value = "'2.085.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").Value = "'4.646"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").Value = "'9.576"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.65%  '

$ws.Range("D24").Value = "'5.997"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").Value = "'144.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").Value = "'2.209"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.42%  '

$ws.Range("D27").Value = "'1.665"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("D28").Value = "'17.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.13%  '

$ws.Range("D29").Value = "'112.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("D30").Value = "'4.211"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("D31").Value = "'4.161"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.30%  '

$ws.Range("D32").Value = "'0.08769"
$ws.Range("D32").ClearFormats()

$ws.Range("D33").Value = "'0.04848"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.77%  '

$ws.Range("D34").Value = "'1.140"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.60%  '

$ws.Range("D35").Value = "'0.7177"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.86%  '

$ws.Range("D36").Value = "'2.855"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("D37").Value = "'3.112"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("D38").Value = "'0.01783"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("D39").Value = "'2.223"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("D40").Value = "'0.4904"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.67%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = "'112.70"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.42%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'0.9075"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("D43").Value = "'6.094"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").Value = "'7.748"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.34%  '

$ws.Range("E46").Value = '  -2.77%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'9.164"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.05942"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("D49").Value = "'0.1244"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.16%  '

$ws.Range("D50").Value = "'35.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("D51").Value = "'0.8901"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.43%  '
